# Natmi LR-pairs (Tnfsf13-Fas) update following Dr Hou advice.
# Replaces the data block (rows 2-17, cols A-T) with recomputed values:
#  - a new "sCs" sending/target cluster category is added (rows 14-17 and the
#    sCs column within each existing sending-cluster block),
#  - expression-count/specificity figures are recalculated accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 16 data rows (sheet rows 2-17) x 20 columns (A-T)
$data = New-Object 'object[,]' 16,20

# sheet row 2
$data[0,0] = "ECs"
$data[0,1] = "Tnfsf13"
$data[0,2] = "Fas"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 0.637903
$data[0,7] = 1.913709
$data[0,8] = 0.1229013127714845
$data[0,9] = 0.1229013127714844
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 18.11506266666667
$data[0,13] = 54.345188
$data[0,14] = 0.5475965948184532
$data[0,15] = 0.5475965948184532
$data[0,16] = 11.55565282025467
$data[0,17] = 104.000875382292
$data[0,18] = 0.06730034037238257
$data[0,19] = 0.06730034037238256

# sheet row 3
$data[1,0] = "ECs"
$data[1,1] = "Tnfsf13"
$data[1,2] = "Fas"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 0.637903
$data[1,7] = 1.913709
$data[1,8] = 0.1229013127714845
$data[1,9] = 0.1229013127714844
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 8.621912666666665
$data[1,13] = 25.865738
$data[1,14] = 0.2606300681353106
$data[1,15] = 0.2606300681353106
$data[1,16] = 5.499943955804666
$data[1,17] = 49.49949560224199
$data[1,18] = 0.03203177752155112
$data[1,19] = 0.03203177752155111

# sheet row 4
$data[2,0] = "ECs"
$data[2,1] = "Tnfsf13"
$data[2,2] = "Fas"
$data[2,3] = "M2"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 0.637903
$data[2,7] = 1.913709
$data[2,8] = 0.1229013127714845
$data[2,9] = 0.1229013127714844
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 3.165745
$data[2,13] = 9.497235
$data[2,14] = 0.09569667044284827
$data[2,15] = 0.09569667044284827
$data[2,16] = 2.019438232735
$data[2,17] = 18.174944094615
$data[2,18] = 0.01176124642528617
$data[2,19] = 0.01176124642528617

# sheet row 5
$data[3,0] = "ECs"
$data[3,1] = "Tnfsf13"
$data[3,2] = "Fas"
$data[3,3] = "sCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 0.637903
$data[3,7] = 1.913709
$data[3,8] = 0.1229013127714845
$data[3,9] = 0.1229013127714844
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 3.178315666666667
$data[3,13] = 9.534947
$data[3,14] = 0.0960766666033877
$data[3,15] = 0.09607666660338772
$data[3,16] = 2.027457098713667
$data[3,17] = 18.247113888423
$data[3,18] = 0.01180794845226459
$data[3,19] = 0.01180794845226459

# sheet row 6
$data[4,0] = "FAPs"
$data[4,1] = "Tnfsf13"
$data[4,2] = "Fas"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 0.170281
$data[4,7] = 0.510843
$data[4,8] = 0.03280711713229307
$data[4,9] = 0.03280711713229307
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 18.11506266666667
$data[4,13] = 54.345188
$data[4,14] = 0.5475965948184532
$data[4,15] = 0.5475965948184532
$data[4,16] = 3.084650985942667
$data[4,17] = 27.761858873484
$data[4,18] = 0.01796506562745382
$data[4,19] = 0.01796506562745382

# sheet row 7
$data[5,0] = "FAPs"
$data[5,1] = "Tnfsf13"
$data[5,2] = "Fas"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 0.170281
$data[5,7] = 0.510843
$data[5,8] = 0.03280711713229307
$data[5,9] = 0.03280711713229307
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 8.621912666666665
$data[5,13] = 25.865738
$data[5,14] = 0.2606300681353106
$data[5,15] = 0.2606300681353106
$data[5,16] = 1.468147910792666
$data[5,17] = 13.213331197134
$data[5,18] = 0.00855052117351266
$data[5,19] = 0.00855052117351266

# sheet row 8
$data[6,0] = "FAPs"
$data[6,1] = "Tnfsf13"
$data[6,2] = "Fas"
$data[6,3] = "M2"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 0.170281
$data[6,7] = 0.510843
$data[6,8] = 0.03280711713229307
$data[6,9] = 0.03280711713229307
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 3.165745
$data[6,13] = 9.497235
$data[6,14] = 0.09569667044284827
$data[6,15] = 0.09569667044284827
$data[6,16] = 0.5390662243450001
$data[6,17] = 4.851596019105
$data[6,18] = 0.003139531876388971
$data[6,19] = 0.003139531876388971

# sheet row 9
$data[7,0] = "FAPs"
$data[7,1] = "Tnfsf13"
$data[7,2] = "Fas"
$data[7,3] = "sCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 0.170281
$data[7,7] = 0.510843
$data[7,8] = 0.03280711713229307
$data[7,9] = 0.03280711713229307
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 3.178315666666667
$data[7,13] = 9.534947
$data[7,14] = 0.0960766666033877
$data[7,15] = 0.09607666660338772
$data[7,16] = 0.5412067700356668
$data[7,17] = 4.870860930321001
$data[7,18] = 0.00315199845493761
$data[7,19] = 0.003151998454937611

# sheet row 10
$data[8,0] = "M2"
$data[8,1] = "Tnfsf13"
$data[8,2] = "Fas"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 3.407124666666667
$data[8,7] = 10.221374
$data[8,8] = 0.6564322386153376
$data[8,9] = 0.6564322386153377
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 18.11506266666667
$data[8,13] = 54.345188
$data[8,14] = 0.5475965948184532
$data[8,15] = 0.5475965948184532
$data[8,16] = 61.72027684981245
$data[8,17] = 555.482491648312
$data[8,18] = 0.3594600585948132
$data[8,19] = 0.3594600585948133

# sheet row 11
$data[9,0] = "M2"
$data[9,1] = "Tnfsf13"
$data[9,2] = "Fas"
$data[9,3] = "FAPs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 3.407124666666667
$data[9,7] = 10.221374
$data[9,8] = 0.6564322386153376
$data[9,9] = 0.6564322386153377
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 8.621912666666665
$data[9,13] = 25.865738
$data[9,14] = 0.2606300681353106
$data[9,15] = 0.2606300681353106
$data[9,16] = 29.37593132044577
$data[9,17] = 264.383381884012
$data[9,18] = 0.1710859790765299
$data[9,19] = 0.17108597907653

# sheet row 12
$data[10,0] = "M2"
$data[10,1] = "Tnfsf13"
$data[10,2] = "Fas"
$data[10,3] = "M2"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 3.407124666666667
$data[10,7] = 10.221374
$data[10,8] = 0.6564322386153376
$data[10,9] = 0.6564322386153377
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 3.165745
$data[10,13] = 9.497235
$data[10,14] = 0.09569667044284827
$data[10,15] = 0.09569667044284827
$data[10,16] = 10.78608787787667
$data[10,17] = 97.07479090089001
$data[10,18] = 0.0628183796068331
$data[10,19] = 0.06281837960683312

# sheet row 13
$data[11,0] = "M2"
$data[11,1] = "Tnfsf13"
$data[11,2] = "Fas"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 3.407124666666667
$data[11,7] = 10.221374
$data[11,8] = 0.6564322386153376
$data[11,9] = 0.6564322386153377
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 3.178315666666667
$data[11,13] = 9.534947
$data[11,14] = 0.0960766666033877
$data[11,15] = 0.09607666660338772
$data[11,16] = 10.82891770635311
$data[11,17] = 97.46025935717802
$data[11,18] = 0.06306782133716124
$data[11,19] = 0.06306782133716125

# sheet row 14
$data[12,0] = "sCs"
$data[12,1] = "Tnfsf13"
$data[12,2] = "Fas"
$data[12,3] = "ECs"
$data[12,4] = 2
$data[12,5] = 0.6666666666666666
$data[12,6] = 0.975059
$data[12,7] = 2.925177
$data[12,8] = 0.1878593314808848
$data[12,9] = 0.1878593314808848
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 18.11506266666667
$data[12,13] = 54.345188
$data[12,14] = 0.5475965948184532
$data[12,15] = 0.5475965948184532
$data[12,16] = 17.66325488869733
$data[12,17] = 158.969293998276
$data[12,18] = 0.1028711302238036
$data[12,19] = 0.1028711302238036

# sheet row 15
$data[13,0] = "sCs"
$data[13,1] = "Tnfsf13"
$data[13,2] = "Fas"
$data[13,3] = "FAPs"
$data[13,4] = 2
$data[13,5] = 0.6666666666666666
$data[13,6] = 0.975059
$data[13,7] = 2.925177
$data[13,8] = 0.1878593314808848
$data[13,9] = 0.1878593314808848
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 8.621912666666665
$data[13,13] = 25.865738
$data[13,14] = 0.2606300681353106
$data[13,15] = 0.2606300681353106
$data[13,16] = 8.406873542847332
$data[13,17] = 75.661861885626
$data[13,18] = 0.04896179036371691
$data[13,19] = 0.04896179036371692

# sheet row 16
$data[14,0] = "sCs"
$data[14,1] = "Tnfsf13"
$data[14,2] = "Fas"
$data[14,3] = "M2"
$data[14,4] = 2
$data[14,5] = 0.6666666666666666
$data[14,6] = 0.975059
$data[14,7] = 2.925177
$data[14,8] = 0.1878593314808848
$data[14,9] = 0.1878593314808848
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 3.165745
$data[14,13] = 9.497235
$data[14,14] = 0.09569667044284827
$data[14,15] = 0.09569667044284827
$data[14,16] = 3.086788153955
$data[14,17] = 27.781093385595
$data[14,18] = 0.01797751253434002
$data[14,19] = 0.01797751253434002

# sheet row 17
$data[15,0] = "sCs"
$data[15,1] = "Tnfsf13"
$data[15,2] = "Fas"
$data[15,3] = "sCs"
$data[15,4] = 2
$data[15,5] = 0.6666666666666666
$data[15,6] = 0.975059
$data[15,7] = 2.925177
$data[15,8] = 0.1878593314808848
$data[15,9] = 0.1878593314808848
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 3.178315666666667
$data[15,13] = 9.534947
$data[15,14] = 0.0960766666033877
$data[15,15] = 0.09607666660338772
$data[15,16] = 3.099045295624333
$data[15,17] = 27.891407660619
$data[15,18] = 0.01804889835902427
$data[15,19] = 0.01804889835902427

# Write the whole block in one shot and refresh the used range.
$ws.Range("A2:T17").Value = $data
